# Applies the two changes described by the diff:
#  1. Bump the cached "datetimeFigureOut" date field text from
#     11.07.2025 -> 14.07.2025 on every slide layout and on the slide
#     master (the Date placeholder shape, ppPlaceholderDate = 16).
#  2. Rename the "Activation" text on slide 12 to "Attention".

$p = $ppt.ActivePresentation

function Get-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        try {
            if ($sh.HasTextFrame -and $sh.PlaceholderFormat.Type -eq 16) {
                return $sh
            }
        } catch {
            # not a placeholder / no PlaceholderFormat - skip
        }
    }
    return $null
}

$oldDate = "11.07.2025"
$newDate = "14.07.2025"

$master = $p.SlideMaster

# 1a. Every slide layout's Date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $dateShape = Get-DatePlaceholder $layout.Shapes
    if ($dateShape -ne $null) {
        if ($dateShape.TextFrame.TextRange.Text -eq $oldDate) {
            $dateShape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# 1b. The slide master's own Date placeholder.
$masterDateShape = Get-DatePlaceholder $master.Shapes
if ($masterDateShape -ne $null) {
    if ($masterDateShape.TextFrame.TextRange.Text -eq $oldDate) {
        $masterDateShape.TextFrame.TextRange.Text = $newDate
    }
}

# 2. Slide 12: "Activation" -> "Attention".
$slide = $p.Slides.Item(12)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "Activation") {
                $sh.TextFrame.TextRange.Text = "Attention"
            }
        }
    }
}
